$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.786.36"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -7.63%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.699.95"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -6.95%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.23%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'569.32"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -6.43%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'172.08"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +2.18%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'3.695.70"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -6.87%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.624"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -8.44%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.997"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.27%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.704"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -10.53%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.164"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -12.05%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'52.35"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -6.54%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.0000295"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -12.38%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'10.52"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -6.79%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'4.304.82"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -6.72%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'3.705.53"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -6.97%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = "'0.127"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -3.03%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = "'19.26"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -6.99%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').Value = "'1.13"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -8.80%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = "'12.84"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -9.98%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'67.734.76"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -7.63%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'404.73"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -11.23%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'4.47"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -7.26%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'87.77"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -8.88%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'3.04"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -11.29%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'12.66"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -11.10%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'10.61"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -4.16%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'3.78"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -9.72%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'5.97"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.07%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'9.45"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -10.21%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'32.61"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -10.44%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'7.57"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -4.26%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'12.52"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -10.05%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.116"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -10.29%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'65.13"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -8.15%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'43.04"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -10.52%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'598.47"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -7.69%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.0₃0887"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -16.16%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = "'0.999"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.05%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = "'0.396"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -8.26%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.15%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.135"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -7.84%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'3.01"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -11.41%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = "'0.0436"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -9.67%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = "'2.86"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -12.76%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'2.56"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.75%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'9.21"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -13.09%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -13.44%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.134"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -10.22%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'2.736.83"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -2.69%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'3.11"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -10.21%  "
$ws.Range('E51').Style = 'Normal'
